$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 1; $r -le 15; $r++) {
    $bText = $ws.Cells.Item($r, 2).Text
    $cText = $ws.Cells.Item($r, 3).Text
    $ws.Cells.Item($r, 2).Value = $cText
    $ws.Cells.Item($r, 3).Value = $bText
}

$ws.Columns.Item(2).ColumnWidth = 29.85546875
$ws.Columns.Item(3).ColumnWidth = 18.140625

$ws.Rows.Item(6).RowHeight = 30

$ws.Range("A10").Select()
$ws.Range("B5").Select()
